$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the two now-unneeded lookup sheets (tblIngredientTypes, tblUnits)
#    - their data has been folded into / is no longer referenced from
#    tblIngredients. Deleting the sheets also removes their associated
#    defined names (tblIngredientTypes, tblUnits) and compacts the shared
#    string table automatically.
$wb.Worksheets("tblIngredientTypes").Delete() | Out-Null
$wb.Worksheets("tblUnits").Delete() | Out-Null

$ws1 = $wb.Worksheets("tblIngredients")
$ws1.Activate()

# 2. Clear the autofilter criteria (column B was filtered down to
#    "Grocery" only) while keeping the filter dropdown arrows on the
#    A1:G48 range. This also unhides the rows that the filter had hidden.
$ws1.ShowAllData()

# 3. Leave the selection/scroll position where the second stage of the
#    expansion (utensils) will continue from.
$ws1.Range("C51").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
